$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 352, shifting existing rows 352:448 down to 353:449
$ws.Rows("352:352").Insert()

# Populate the newly inserted row 352 with its data
$ws.Range("A352").Value = 5
$ws.Range("B352").Value = "Macroferia Regional de Talca"
$ws.Range("C352").Value = "Maule"
$ws.Range("D352").Value = 44551
$ws.Range("E352").Value = 7
$ws.Range("F352").Value = "Fruta"
$ws.Range("G352").Value = 100104
$ws.Range("H352").Value = "Frutos de pepita"
$ws.Range("I352").Value = 100104005
$ws.Range("J352").Value = "Pera"
$ws.Range("K352").Value = "Packham's Triumph"
$ws.Range("L352").Value = "Primera"
$ws.Range("M352").Value = 230
$ws.Range("N352").Value = 9000
$ws.Range("O352").Value = 9000
$ws.Range("P352").Value = 9000
$ws.Range("Q352").Value = "$/bandeja 18 kilos granel"
$ws.Range("R352").Value = "Provincia de Curicó"
$ws.Range("S352").Value = 500
$ws.Range("T352").Value = 18
